# Auto-generated edit applying the cryptos.xlsx price/volume refresh diff.
# D-column values are prefixed with a leading apostrophe so Excel keeps them
# as text (matching the original inline-string cells) instead of auto-converting
# number-looking strings (e.g. "321.30", "26.00") into numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'46.291.52"
$ws.Range("E2").Value = "  +3.80%  "

$ws.Range("D3").Value = "'2.454.71"
$ws.Range("E3").Value = "  +1.30%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'321.30"
$ws.Range("E5").Value = "  +1.90%  "

$ws.Range("D6").Value = "'105.54"
$ws.Range("E6").Value = "  +4.40%  "

$ws.Range("E7").Value = "  +0.78%  "

$ws.Range("D9").Value = "'0.536"
$ws.Range("E9").Value = "  +2.19%  "

$ws.Range("D10").Value = "'36.10"
$ws.Range("E10").Value = "  +2.07%  "

$ws.Range("D11").Value = "'0.0815"
$ws.Range("E11").Value = "  +1.88%  "

$ws.Range("D12").Value = "'0.123"
$ws.Range("E12").Value = "  +0.97%  "

$ws.Range("D13").Value = "'18.33"
$ws.Range("E13").Value = "  -3.68%  "

$ws.Range("D14").Value = "'7.07"
$ws.Range("E14").Value = "  +1.34%  "

$ws.Range("D15").Value = "'2.844.28"
$ws.Range("E15").Value = "  +1.54%  "

$ws.Range("D16").Value = "'2.430.33"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("E17").Value = "  +1.35%  "

$ws.Range("D18").Value = "'46.125.72"
$ws.Range("E18").Value = "  +3.73%  "

$ws.Range("D19").Value = "'12.75"
$ws.Range("E19").Value = "  +2.23%  "

$ws.Range("E20").Value = "  +0.44%  "

$ws.Range("D21").Value = "'0.0₃0937"
$ws.Range("E21").Value = "  +1.86%  "

$ws.Range("D22").Value = "'70.52"
$ws.Range("E22").Value = "  +2.56%  "

$ws.Range("E23").Value = "  +4.41%  "

$ws.Range("D24").Value = "'247.91"
$ws.Range("E24").Value = "  +2.21%  "

$ws.Range("D25").Value = "'2.53"
$ws.Range("E25").Value = "  +2.07%  "

$ws.Range("D26").Value = "'26.00"
$ws.Range("E26").Value = "  +2.94%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("E28").Value = "  +0.49%  "

$ws.Range("E29").Value = "  +1.66%  "

$ws.Range("D30").Value = "'34.76"
$ws.Range("E30").Value = "  +4.83%  "

$ws.Range("D31").Value = "'49.35"
$ws.Range("E31").Value = "  +1.84%  "

$ws.Range("E32").Value = "  +3.74%  "

$ws.Range("D33").Value = "'19.78"
$ws.Range("E33").Value = "  +2.40%  "

$ws.Range("E34").Value = "  +3.30%  "

$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").Value = "'0.0767"
$ws.Range("E36").Value = "  -1.31%  "

$ws.Range("D37").Value = "'4.58"
$ws.Range("E37").Value = "  +1.84%  "

$ws.Range("E38").Value = "  +0.54%  "

$ws.Range("D39").Value = "'2.96"
$ws.Range("E39").Value = "  +3.25%  "

$ws.Range("D40").Value = "'124.23"
$ws.Range("E40").Value = "  +4.48%  "

$ws.Range("E41").Value = "  +1.58%  "

$ws.Range("D42").Value = "'2.23"
$ws.Range("E42").Value = "  +0.23%  "

$ws.Range("D43").Value = "'20.79"
$ws.Range("E43").Value = "  -0.62%  "

$ws.Range("E44").Value = "  +0.70%  "

$ws.Range("D45").Value = "'1.976.86"
$ws.Range("E45").Value = "  +1.85%  "

$ws.Range("E46").Value = "  +1.27%  "

$ws.Range("E47").Value = "  -3.77%  "

$ws.Range("E48").Value = "  +10.92%  "

$ws.Range("D49").Value = "'9.09"
$ws.Range("E49").Value = "  -3.54%  "

$ws.Range("D50").Value = "'5.11"
$ws.Range("E50").Value = "  +10.31%  "

$ws.Range("D51").Value = "'78.59"
$ws.Range("E51").Value = "  +4.90%  "
